# Auto-generated edit script applying the cryptos.xlsx data refresh diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.107.24'
$ws.Range('E2').Value = '  -1.90%  '
$ws.Range('D3').Value = '2.575.51'
$ws.Range('E3').Value = '  -3.05%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '585.76'
$ws.Range('E5').Value = '  -3.62%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '147.82'
$ws.Range('E6').Value = '  -3.38%  '
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('E8').Value = '  -1.31%  '
$ws.Range('E9').Value = '  -1.00%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.71'
$ws.Range('E10').Value = '  +1.90%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.380'
$ws.Range('E11').Value = '  -1.64%  '
$ws.Range('E12').Value = '  -0.85%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '27.27'
$ws.Range('E13').Value = '  -3.35%  '
$ws.Range('D14').Value = '3.036.26'
$ws.Range('E14').Value = '  -3.16%  '
$ws.Range('D15').Value = '62.996.22'
$ws.Range('E15').Value = '  -2.03%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000152'
$ws.Range('E16').Value = '  +2.63%  '
$ws.Range('D17').Value = '2.563.40'
$ws.Range('E17').Value = '  -3.63%  '
$ws.Range('E18').Value = '  -0.31%  '
$ws.Range('E19').Value = '  +0.04%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '343.53'
$ws.Range('E20').Value = '  -1.57%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.79'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.998'
$ws.Range('E22').Value = '  -0.09%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '66.71'
$ws.Range('E23').Value = '  -0.17%  '
$ws.Range('E24').Value = '  -3.08%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.04'
$ws.Range('E25').Value = '  -3.79%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.63'
$ws.Range('E26').Value = '  -4.40%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '552.35'
$ws.Range('E27').Value = '  -0.14%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.99'
$ws.Range('E28').Value = '  -2.48%  '
$ws.Range('B29').Value = 'Binance-PegBSC-USD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  +0.18%  '
$ws.Range('B30').Value = 'Kaspa'
$ws.Range('C30').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.161'
$ws.Range('E30').Value = '  -2.29%  '
$ws.Range('E31').Value = '  -2.60%  '
$ws.Range('D32').Value = '0.0₃0847'
$ws.Range('E32').Value = '  -1.97%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.73'
$ws.Range('E33').Value = '  -2.20%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.14'
$ws.Range('E34').Value = '  -4.30%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '165.10'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.410'
$ws.Range('E36').Value = '  +0.34%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.999'
$ws.Range('E37').Value = '  -0.11%  '
$ws.Range('E38').Value = '  -0.51%  '
$ws.Range('E39').Value = '  -4.53%  '
$ws.Range('E40').Value = '  -0.01%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '165.28'
$ws.Range('E41').Value = '  -0.77%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '39.54'
$ws.Range('E42').Value = '  -1.56%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.94'
$ws.Range('E43').Value = '  +2.11%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0584'
$ws.Range('E44').Value = '  +0.90%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '22.42'
$ws.Range('E45').Value = '  +1.87%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.627'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.02'
$ws.Range('E47').Value = '  +0.94%  '
$ws.Range('E48').Value = '  -0.33%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0958'
$ws.Range('E49').Value = '  -0.98%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '18.86'
$ws.Range('E50').Value = '  -1.24%  '
$ws.Range('E51').Value = '  +10.89%  '
